$d = $word.ActiveDocument
$checkbox = [char]0x2610

# ---------------------------------------------------------------------
# Change 1: the first "Predicted Voltages..." checklist item (numId 1004)
# becomes the calibration-curve item, and a brand-new checklist item
# ("Predicted and measured Voltages match...") is inserted right after it.
# (There is a second, untouched, copy of the original sentence later in
# the doc under numId 1007 - use wdReplaceOne starting from the top so
# only the first hit is touched.)
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "Predicted Voltages and measured voltages match to within uncertainty and are consistent with equations.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $r1.Text = "A graph of a calibration curve between pressure and depth is provided. The data points match the points in the table below. Appropriate regressions and measures of uncertainty are plotted."
    $p1 = $r1.Paragraphs(1)
    $p1.Range.InsertParagraphAfter()
    $p1new = $p1.Next()
    $p1new.Range.Text = $checkbox + " Predicted and measured Voltages match to within uncertainty and are consistent with equations."
}

# ---------------------------------------------------------------------
# Change 2: swap "including the pressure sensor" for
# "including the thermistor divider." on the thermistor-circuit-picture
# checklist item (unique occurrence, so a plain replace-all is safe).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Provides a picture of the thermistor interface circuit, including the pressure sensor",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Provides a picture of the thermistor interface circuit, including the thermistor divider.", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3: insert a new checklist item right after "Settling times are
# consistent with results from lecture." (unique occurrence).
# ---------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "Settling times are consistent with results from lecture.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found3) {
    $p3 = $r3.Paragraphs(1)
    $p3.Range.InsertParagraphAfter()
    $p3new = $p3.Next()
    $p3new.Range.Text = $checkbox + " Brief (3-5 sentence) discussion about which sensors are appropriate for E80 robots and for cold-junction compensation sensors."
}
